$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 286
$ws.Range("G8").Value = 98
$ws.Range("F9").Value = 125
$ws.Range("F11").Value = 5181
$ws.Range("F13").Value = 3621
$ws.Range("F14").Value = 53
$ws.Range("F17").Value = 182
$ws.Range("F21").Value = 80
$ws.Range("F22").Value = 131
$ws.Range("F23").Value = 90
$ws.Range("F25").Value = 4694
$ws.Range("F27").Value = 1973
$ws.Range("F28").Value = 113
$ws.Range("F30").Value = 7178
$ws.Range("F33").Value = 2149
$ws.Range("F34").Value = 2072
$ws.Range("F35").Value = 1311
$ws.Range("F36").Value = 125
$ws.Range("F37").Value = 1116
$ws.Range("F38").Value = 17
$ws.Range("F39").Value = 8
$ws.Range("F43").Value = 12
$ws.Range("F44").Value = 156
$ws.Range("F45").Value = 1258
$ws.Range("F46").Value = 1887
$ws.Range("F47").Value = 84
$ws.Range("F49").Value = 1191

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 496
$ws.Range("F3").Value = 671

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 496
$ws.Range("F6").Value = 671
$ws.Range("F7").Value = 286
$ws.Range("F11").Value = 5181
$ws.Range("F12").Value = 3621
$ws.Range("F13").Value = 53
$ws.Range("F15").Value = 182
$ws.Range("F18").Value = 80
$ws.Range("F20").Value = 131
$ws.Range("F21").Value = 90
$ws.Range("F24").Value = 4694
$ws.Range("F26").Value = 1973
$ws.Range("F27").Value = 113
$ws.Range("F29").Value = 7178
$ws.Range("F33").Value = 2149
$ws.Range("F34").Value = 2072
$ws.Range("F35").Value = 1311
$ws.Range("F36").Value = 125
$ws.Range("F37").Value = 1116
$ws.Range("F38").Value = 8
$ws.Range("F42").Value = 156
$ws.Range("F44").Value = 1258
$ws.Range("F45").Value = 1887
$ws.Range("F46").Value = 84
$ws.Range("F49").Value = 1191

Write-Output "Applied F/G count updates across sheets 1, 3, 4."
